$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (row 1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 23:29"

# Refresh country statistics. Some countries reorder rows as a side effect
# of the data refresh (the sheet is kept sorted by "Casos totales" descending),
# so both the country name (column A) and its stats (B:H) are written per row.
$countryData = @(
    @(4, "Estados Unidos", 4962826, 44406, 2524395, 2277029, 0, 1112, 161402),
    @(5, "Brasil", 2858872, 50796, 1970767, 790817, 0, 1192, 97288),
    @(8, "Sudafrica", 529877, 8559, 377266, 143313, 0, 414, 9298),
    @(30, "Egipto", 94875, 123, 47182, 42763, 0, 18, 4930),
    @(35, "Oman", 80286, 585, 69803, 9995, 0, 6, 488),
    @(58, "Suiza", 35927, 181, 31600, 2343, 0, 3, 1984),
    @(134, "Benin", 1914, 0, 1600, 276, 0, 0, 38),
    @(136, "Yemen", 1763, 3, 894, 361, 0, 2, 508),
    @(146, "Republica de Chipre", 1195, 15, 856, 320, 0, 0, 19),
    @(149, "Togo", 1001, 13, 690, 290, 0, 2, 21),
    @(152, "Republica del Chad", 939, 1, 835, 29, 0, 0, 75),
    @(155, "Santo Tome y Principe", 878, 3, 795, 68, 0, 0, 15),
    @(183, "San Martin (Parte Holandesa)", 156, 6, 64, 76, 0, 0, 16),
    @(184, "Papua Nueva Guinea", 153, 39, 44, 107, 0, 0, 2),
    @(186, "Barbados", 133, 1, 100, 26, 0, 0, 7),
    @(187, "Seychelles", 126, 12, 124, 2, 0, 0, 0),
    @(188, "Monaco", 125, 2, 105, 16, 0, 0, 4),
    @(189, "Islas Turcas y Caicos", 116, 0, 38, 76, 0, 0, 2),
    @(202, "Timor Oriental", 25, 0, 24, 1, 0, 0, 0),
    @(203, "Santa Lucia", 25, 0, 24, 1, 0, 0, 0)
)

foreach ($entry in $countryData) {
    $rowNum = $entry[0]
    $ws.Cells.Item($rowNum, 1).Value = $entry[1]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $entry[$col]
    }
}
